# Auto-generated from the cryptos.xlsx price/volume refresh diff.
# Column D = Price, column E = Volume(1h); both are plain text cells
# (note the tell-tale "63.933.00"-style thousands separators and the
# padded "  +1.41%  " volume strings -- neither parses as a normal number).
#
# Excel auto-detects numeric-looking text on assignment (e.g. "34.80" or
# "0.0000249") and would silently coerce it to a float, corrupting the
# trailing/leading zeros. We guard those cells with a leading apostrophe
# (the same text-prefix a user would type in the Excel UI) so the value is
# stored verbatim as text. Values that already fail numeric parsing (e.g.
# "63.790.52", which has two dots) do not need the guard.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.790.52'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.313.56'
$ws.Range("E3").Value = '  +6.26%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''601.22'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").Value = '''142.87'
$ws.Range("E6").Value = '  +4.89%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '3.308.66'
$ws.Range("E8").Value = '  +6.30%  '
$ws.Range("D9").Value = '''0.522'
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").Value = '''0.149'
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("D11").Value = '''5.53'
$ws.Range("E11").Value = '  +4.36%  '
$ws.Range("D12").Value = '''0.473'
$ws.Range("E12").Value = '  +4.22%  '
$ws.Range("D13").Value = '''0.0000249'
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("D14").Value = '''34.80'
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '3.866.01'
$ws.Range("E15").Value = '  +6.45%  '
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").Value = '3.315.16'
$ws.Range("E17").Value = '  +6.14%  '
$ws.Range("D18").Value = '64.019.70'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '''6.89'
$ws.Range("E19").Value = '  +3.65%  '
$ws.Range("D20").Value = '''480.95'
$ws.Range("E20").Value = '  +2.00%  '
$ws.Range("D21").Value = '''14.23'
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").Value = '''0.736'
$ws.Range("E22").Value = '  +6.00%  '
$ws.Range("D23").Value = '''7.98'
$ws.Range("E23").Value = '  +3.70%  '
$ws.Range("D24").Value = '''13.54'
$ws.Range("E24").Value = '  +5.31%  '
$ws.Range("D25").Value = '''84.58'
$ws.Range("E25").Value = '  -1.81%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '''2.77'
$ws.Range("E27").Value = '  +2.12%  '
$ws.Range("D28").Value = '''7.33'
$ws.Range("E28").Value = '  +5.73%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = '''8.17'
$ws.Range("E30").Value = '  +3.51%  '
$ws.Range("D31").Value = '''2.16'
$ws.Range("E31").Value = '  +5.01%  '
$ws.Range("D32").Value = '''29.26'
$ws.Range("E32").Value = '  +10.14%  '
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("D34").Value = '''2.56'
$ws.Range("E34").Value = '  +1.72%  '
$ws.Range("D35").Value = '''1.10'
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("D36").Value = '''5.98'
$ws.Range("E36").Value = '  +3.56%  '
$ws.Range("D37").Value = '0.0₃0753'
$ws.Range("E37").Value = '  +8.19%  '
$ws.Range("D38").Value = '''52.76'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = '''0.0405'
$ws.Range("E39").Value = '  +5.16%  '
$ws.Range("D40").Value = '''431.38'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '3.055.97'
$ws.Range("E41").Value = '  +5.71%  '
$ws.Range("D42").Value = '''8.42'
$ws.Range("E42").Value = '  +3.04%  '
$ws.Range("D43").Value = '''2.76'
$ws.Range("E43").Value = '  +3.16%  '
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '''0.267'
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("D46").Value = '''2.20'
$ws.Range("E46").Value = '  +4.65%  '
$ws.Range("D47").Value = '''26.45'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("D48").Value = '''36.05'
$ws.Range("E48").Value = '  +15.12%  '
$ws.Range("E50").Value = '  +2.20%  '
$ws.Range("D51").Value = '''2.30'
$ws.Range("E51").Value = '  +2.33%  '
